$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 18 should look like row 14 (style-wise: s="6" on A/B/D/E, s="2" on C)
$ws.Range("A14:E14").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)

$ws.Range("A18").Value = "TestCase_F17"
$ws.Range("B18").Value = "OPQA-1098"
$ws.Range("C18").Value = "Verify that Featured Post is at the top of event stream after login and that feature post should be top in post tab of trending section"
$ws.Range("D18").Value = "Y"
$ws.Range("E18").Value = "PASS"

$ws.Range("D17").Select()
